$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 value from 2 to 3
$ws.Range("B2").Value = 3

# Remove row 3 entirely (A3=1, B3=1 -> deleted), shrinking the used range to A1:B2
$ws.Rows("3").Delete()
